$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $ws.Range("ZZ1").NumberFormat = "@"
    $ws.Range("ZZ1").Value = $val
    $ws.Range("ZZ1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $ws.Range("ZZ1").Clear()
}

# Row 2
Set-TextValue 'D2' '22.039.98'
Set-TextValue 'E2' '  -1.54%  '

# Row 3
Set-TextValue 'D3' '1.554.66'
Set-TextValue 'E3' '  -0.79%  '

# Row 4
Set-TextValue 'E4' '  -0.05%  '

# Row 5
Set-TextValue 'E5' '  +0.02%  '

# Row 6
Set-TextValue 'D6' '287.95'
Set-TextValue 'E6' '  +0.48%  '

# Row 7
Set-TextValue 'D7' '0.3984'
Set-TextValue 'E7' '  +6.43%  '

# Row 8
Set-TextValue 'D8' '0.3224'
Set-TextValue 'E8' '  -1.53%  '

# Row 9
Set-TextValue 'D9' '42.54'
Set-TextValue 'E9' '  -6.31%  '

# Row 10
Set-TextValue 'B10' 'Dogecoin'
Set-TextValue 'C10' 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue 'D10' '0.07337'
Set-TextValue 'E10' '  -1.08%  '

# Row 11
Set-TextValue 'B11' 'Polygon'
Set-TextValue 'C11' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D11' '1.106'
Set-TextValue 'E11' '  -3.94%  '

# Row 12
Set-TextValue 'E12' '  -0.06%  '

# Row 13
Set-TextValue 'D13' '18.95'

# Row 14
Set-TextValue 'D14' '5.684'
Set-TextValue 'E14' '  -2.68%  '

# Row 15
Set-TextValue 'D15' '6.736'
Set-TextValue 'E15' '  -1.47%  '

# Row 16
Set-TextValue 'B16' 'ShibaInu'
Set-TextValue 'C16' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D16' '0.00001130'
Set-TextValue 'E16' '  +2.86%  '

# Row 17
Set-TextValue 'B17' 'WrappedEther'
Set-TextValue 'C17' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D17' '1.552.97'
Set-TextValue 'E17' '  -0.79%  '

# Row 18
Set-TextValue 'D18' '0.06608'
Set-TextValue 'E18' '  -1.56%  '

# Row 19
Set-TextValue 'D19' '85.36'
Set-TextValue 'E19' '  -0.58%  '

# Row 20
Set-TextValue 'E20' '  +0.01%  '

# Row 21
Set-TextValue 'D21' '6.352'
Set-TextValue 'E21' '  +0.06%  '

# Row 22
Set-TextValue 'E22' '  -2.21%  '

# Row 23
Set-TextValue 'D23' '11.29'
Set-TextValue 'E23' '  -3.48%  '

# Row 24
Set-TextValue 'D24' '22.049.48'
Set-TextValue 'E24' '  -1.52%  '

# Row 25
Set-TextValue 'D25' '2.331'
Set-TextValue 'E25' '  +0.77%  '

# Row 26
Set-TextValue 'D26' '2.478'
Set-TextValue 'E26' '  -2.93%  '

# Row 27
Set-TextValue 'D27' '148.21'
Set-TextValue 'E27' '  -2.02%  '

# Row 28
Set-TextValue 'D28' '18.75'
Set-TextValue 'E28' '  -3.12%  '

# Row 29
Set-TextValue 'D29' '4.862'
Set-TextValue 'E29' '  -1.04%  '

# Row 30
Set-TextValue 'D30' '1.732.37'
Set-TextValue 'E30' '  -0.63%  '

# Row 31
Set-TextValue 'D31' '120.75'
Set-TextValue 'E31' '  -2.13%  '

# Row 32
Set-TextValue 'D32' '1.071'
Set-TextValue 'E32' '  +1.39%  '

# Row 33
Set-TextValue 'D33' '5.728'
Set-TextValue 'E33' '  -3.29%  '

# Row 34
Set-TextValue 'D34' '0.08409'
Set-TextValue 'E34' '  +1.77%  '

# Row 35
Set-TextValue 'B35' 'FraxShare'
Set-TextValue 'C35' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D35' '9.318'
Set-TextValue 'E35' '  -3.26%  '

# Row 36
Set-TextValue 'B36' 'WEMIXTOKEN'
Set-TextValue 'C36' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D36' '1.640'
Set-TextValue 'E36' '  -15.57%  '

# Row 37
Set-TextValue 'D37' '0.06231'
Set-TextValue 'E37' '  -1.46%  '

# Row 38
Set-TextValue 'D38' '0.02274'
Set-TextValue 'E38' '  -4.71%  '

# Row 39
Set-TextValue 'D39' '5.144'
Set-TextValue 'E39' '  -2.06%  '

# Row 40
Set-TextValue 'D40' '0.2095'
Set-TextValue 'E40' '  -4.23%  '

# Row 41
Set-TextValue 'D41' '1.214'
Set-TextValue 'E41' '  -5.93%  '

# Row 42
Set-TextValue 'D42' '1.000'
Set-TextValue 'E42' '  -0.01%  '

# Row 43
Set-TextValue 'D43' '10.67'
Set-TextValue 'E43' '  -3.88%  '

# Row 44
Set-TextValue 'D44' '0.5878'
Set-TextValue 'E44' '  -3.66%  '

# Row 45
Set-TextValue 'D45' '13.32'
Set-TextValue 'E45' '  -2.78%  '

# Row 46
Set-TextValue 'D46' '3.719'
Set-TextValue 'E46' '  -0.75%  '

# Row 47
Set-TextValue 'D47' '0.5607'
Set-TextValue 'E47' '  -5.26%  '

# Row 48
Set-TextValue 'D48' '1.914'
Set-TextValue 'E48' '  -4.71%  '

# Row 49
Set-TextValue 'D49' '117.92'
Set-TextValue 'E49' '  -4.74%  '

# Row 50
Set-TextValue 'D50' '1.148'
Set-TextValue 'E50' '  -2.72%  '

# Row 51
Set-TextValue 'E51' '  -4.18%  '
